$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.238.25"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").Value = "1.649.68"
$ws.Range("E3").Value = "  -0.01%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'218.77"
$ws.Range("E5").Value = "  -0.42%  "
$ws.Range("E6").Value = "  +2.08%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'0.257"
$ws.Range("E8").Value = "  +0.95%  "
$ws.Range("E9").Value = "  +0.35%  "
$ws.Range("D10").Value = "'20.24"
$ws.Range("E10").Value = "  +2.77%  "
$ws.Range("E11").Value = "  -0.04%  "
$ws.Range("D12").Value = "1.880.78"
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("D13").Value = "1.647.69"
$ws.Range("E13").Value = "  -0.09%  "
$ws.Range("D14").Value = "'4.14"
$ws.Range("E14").Value = "  -1.35%  "
$ws.Range("D15").Value = "'0.539"
$ws.Range("E15").Value = "  +0.94%  "
$ws.Range("D16").Value = "'68.09"
$ws.Range("E16").Value = "  +2.74%  "
$ws.Range("D17").Value = "27.210.05"
$ws.Range("E17").Value = "  +0.39%  "
$ws.Range("E18").Value = "  +0.46%  "
$ws.Range("D19").Value = "'221.08"
$ws.Range("E19").Value = "  -0.33%  "
$ws.Range("E20").Value = "  +0.08%  "
$ws.Range("D21").Value = "'6.75"
$ws.Range("E21").Value = "  -0.83%  "
$ws.Range("D23").Value = "'2.50"
$ws.Range("E23").Value = "  +3.40%  "
$ws.Range("E24").Value = "  -0.25%  "
$ws.Range("D25").Value = "'148.26"
$ws.Range("E25").Value = "  +0.58%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("E28").Value = "  +0.60%  "
$ws.Range("D29").Value = "'15.84"
$ws.Range("E29").Value = "  -0.44%  "
$ws.Range("E30").Value = "  -1.07%  "
$ws.Range("E31").Value = "  -0.39%  "
$ws.Range("E33").Value = "  +0.35%  "
$ws.Range("E34").Value = "  +0.52%  "
$ws.Range("D35").Value = "1.270.34"
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("D36").Value = "'2.46"
$ws.Range("E36").Value = "  +1.21%  "
$ws.Range("E37").Value = "  +2.32%  "
$ws.Range("D38").Value = "'0.545"
$ws.Range("E38").Value = "  +1.21%  "
$ws.Range("D39").Value = "'0.845"
$ws.Range("E39").Value = "  +2.15%  "
$ws.Range("E40").Value = "  +0.07%  "
$ws.Range("E41").Value = "  +0.38%  "
$ws.Range("E42").Value = "  +0.51%  "
$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D43").Value = "'2.19"
$ws.Range("E43").Value = "  +6.14%  "
$ws.Range("B44").Value = "RocketPoolETH"
$ws.Range("C44").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D44").Value = "1.791.81"
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("D45").Value = "'63.14"
$ws.Range("E45").Value = "  +1.81%  "
$ws.Range("D46").Value = "'92.59"
$ws.Range("E46").Value = "  -0.05%  "
$ws.Range("E47").Value = "  -0.93%  "
$ws.Range("E48").Value = "  +16.77%  "
$ws.Range("D49").Value = "'0.0513"
$ws.Range("E49").Value = "  -0.54%  "
$ws.Range("D50").Value = "'7.72"
$ws.Range("E50").Value = "  +1.11%  "
$ws.Range("D51").Value = "'0.0977"
$ws.Range("E51").Value = "  +0.23%  "
